$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Plane" values for the B-mics (rows 2-18) and D-mics (rows 36-52)
# were swapped: B-mics were labeled "right" but should be "left", and
# D-mics were labeled "left" but should be "right".
$ws.Range("C2:C18").Value = "left"
$ws.Range("C36:C52").Value = "right"
